# DALmethodes.xlsx - "DAL download en subfamily gemaakt + getest" edit
#
# Fills in the Auteur/Status columns for the DaDownload (rows 25-30) and
# DaSubFamily (rows 37-42) sections with "Oualid" / "Done", and updates the
# DaUserAccount.selectOneByUsername() author cell (C60) from "lenny" to
# "lenny/Oualid" to reflect pair-programming / shared credit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- C60: lenny -> lenny/Oualid (do this first so the new shared string is
#     appended before "Oualid", matching the authoring order of the commit) ---
$ws.Range("C60").Value = "lenny/Oualid"

# --- DaDownload rows (25-30): Auteur = Oualid, Status = Done ---
$downloadRows = 25,26,27,28,29,30
foreach ($r in $downloadRows) {
    $ws.Cells.Item($r, 3).Value = "Oualid"
    $ws.Cells.Item($r, 4).Value = "Done"
}

# --- DaSubFamily rows (37-42): Auteur = Oualid, Status = Done ---
$subfamilyRows = 37,38,39,40,41,42
foreach ($r in $subfamilyRows) {
    $ws.Cells.Item($r, 3).Value = "Oualid"
    $ws.Cells.Item($r, 4).Value = "Done"
}

# --- Formatting: the new cells pick up the same thin/no border + default
#     (non-applied) font treatment used elsewhere in the sheet for freshly
#     filled-in status cells, distinguishing them from the plain body style. ---
$newCells = $ws.Range("C25:D30,C37:D42,C60")
$newCells.BorderAround(-4142)

# --- View state: scrolled down to the newly completed sections, selection
#     left on the last touched cell ---
$ws.Application.ActiveWindow.ScrollRow = 28
$ws.Range("D40").Select()
